# Fixed bug in Vect operator+= and added VolPyra and VolPrism
# The volume calcs are based on a degenerate hex volume. Adds a new
# "Sheet2" worksheet (after "Sheet1") containing the volume calculations,
# and makes it the active/selected sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Add the new worksheet right after Sheet1 and name it "Sheet2".
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# --- Row 8: base vertex (all zero) ---
$ws2.Range("B8").Value = 0
$ws2.Range("C8").Value = 0
$ws2.Range("D8").Value = 0
$ws2.Range("E8").Value = 0

# --- Row 9 ---
$ws2.Range("B9").Formula = "=B8+1"
$ws2.Range("C9").Value = 2
$ws2.Range("D9").Value = 0
$ws2.Range("E9").Value = 0

# --- Row 10 ---
$ws2.Range("B10").Formula = "=B9+1"
$ws2.Range("C10").Value = 2
$ws2.Range("D10").Value = 1
$ws2.Range("E10").Value = 0

# --- Row 11 ---
$ws2.Range("B11").Formula = "=B10+1"
$ws2.Range("C11").Value = 0
$ws2.Range("D11").Value = 1
$ws2.Range("E11").Value = 0

# --- Row 12 ---
$ws2.Range("B12").Formula = "=B11+1"
$ws2.Range("C12").Value = 0
$ws2.Range("D12").Value = 0
$ws2.Range("E12").Value = 3

# --- Row 13 ---
$ws2.Range("B13").Formula = "=B12+1"
$ws2.Range("C13").Value = 2
$ws2.Range("D13").Value = 0
$ws2.Range("E13").Value = 3

# --- Row 14 ---
$ws2.Range("B14").Formula = "=B13+1"
$ws2.Range("C14").Value = 2
$ws2.Range("D14").Value = 1
$ws2.Range("E14").Value = 3

# --- Row 15 ---
$ws2.Range("B15").Formula = "=B14+1"
$ws2.Range("C15").Value = 0
$ws2.Range("D15").Value = 1
$ws2.Range("E15").Value = 3

# --- Row 18: diagonal 1 vector ---
$ws2.Range("C18").Formula = "=C14-C9+C15-C8"
$ws2.Range("D18").Formula = "=D14-D9+D15-D8"
$ws2.Range("E18").Formula = "=E14-E9+E15-E8"

# --- Row 19: labelled "a" ---
$ws2.Range("B19").Value = "a"
$ws2.Range("C19").Formula = "=C14-C11"
$ws2.Range("D19").Formula = "=D14-D11"
$ws2.Range("E19").Formula = "=E14-E11"

# --- Row 20: labelled "b" (right-aligned numbers) ---
$ws2.Range("B20").Value = "b"
$ws2.Range("C20").Formula = "=C10-C8"
$ws2.Range("D20").Formula = "=D10-D8"
$ws2.Range("E20").Formula = "=E10-E8"
$ws2.Range("C20:E20").HorizontalAlignment = -4152

# --- Row 22: triple product (tetra volume piece 1) ---
$ws2.Range("C22").Formula = "=(D19*E20-E19*D20)*C18+(E19*C20-C19*E20)*D18+(C19*D20-D19*C20)*E18"

# --- Row 24-26 ---
$ws2.Range("C24").Formula = "=C15-C8"
$ws2.Range("D24").Formula = "=D15-D8"
$ws2.Range("E24").Formula = "=E15-E8"

$ws2.Range("C25").Formula = "=C14-C11+C13-C8"
$ws2.Range("D25").Formula = "=D14-D11+D13-D8"
$ws2.Range("E25").Formula = "=E14-E11+E13-E8"

$ws2.Range("C26").Formula = "=C14-C12"
$ws2.Range("D26").Formula = "=D14-D12"
$ws2.Range("E26").Formula = "=E14-E12"

# --- Row 28: tetra volume piece 2 ---
$ws2.Range("C28").Formula = "=(D25*E26-D26*E25)*C24-(C25*E26-C26*E25)*D24+(C25*D26-C26*D25)*E24"

# --- Row 30-32 ---
$ws2.Range("C30").Formula = "=C14-C9"
$ws2.Range("D30").Formula = "=D14-D9"
$ws2.Range("E30").Formula = "=E14-E9"

$ws2.Range("C31").Formula = "=C13-C8"
$ws2.Range("D31").Formula = "=D13-D8"
$ws2.Range("E31").Formula = "=E13-E8"

$ws2.Range("C32").Formula = "=C14-C12+C10-C8"
$ws2.Range("D32").Formula = "=D14-D12+D10-D8"
$ws2.Range("E32").Formula = "=E14-E12+E10-E8"

# --- Row 34: tetra volume piece 3 ---
$ws2.Range("C34").Formula = "=(D31*E32-D32*E31)*C30-(C31*E32-C32*E31)*D30+(C31*D32-C32*D31)*E30"

# --- Row 36-37: sum and divide by 12 to get the hex volume ---
$ws2.Range("C36").Formula = "=C22+C28+C34"
$ws2.Range("C37").Formula = "=C36/12"

# Page setup (best-effort cosmetic match to the authored sheet view / print
# settings; none of this affects any computed cell values).
$ps = $ws2.PageSetup
$ps.LeftMargin = 56.7
$ps.RightMargin = 56.7
$ps.TopMargin = 75.8
$ps.BottomMargin = 75.8
$ps.HeaderMargin = 56.7
$ps.FooterMargin = 56.7
$ps.Zoom = 100
$ps.CenterHeader = '&"Times New Roman,Regular"&12&A'
$ps.CenterFooter = '&"Times New Roman,Regular"&12Page &P'

# Selection / active-sheet bookkeeping to match the authored view state.
$ws2.Range("B5").Select()
$ws2.Activate()
